$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column F - copy formatting (bold header style) from E1, then set the text
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# time_taken values for rows 2-18
$times = @(
    "2021-10-05 13:41:16.197187",
    "2021-10-05 13:41:16.197197",
    "2021-10-05 13:41:16.197201",
    "2021-10-05 13:41:16.197203",
    "2021-10-05 13:41:16.197206",
    "2021-10-05 13:41:16.197209",
    "2021-10-05 13:41:16.197212",
    "2021-10-05 13:41:16.197214",
    "2021-10-05 13:41:16.197217",
    "2021-10-05 13:41:16.197220",
    "2021-10-05 13:41:16.197222",
    "2021-10-05 13:41:16.197225",
    "2021-10-05 13:41:16.197228",
    "2021-10-05 13:41:16.197230",
    "2021-10-05 13:41:16.197233",
    "2021-10-05 13:41:16.197235",
    "2021-10-05 13:41:16.197238"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
